# Updates the cryptos price table (columns B-E, rows 2-51) with refreshed
# values. Each assigned string is prefixed with a leading apostrophe
# (escaped as '' inside the single-quoted PowerShell literal) so Excel
# stores it as text instead of auto-converting number-looking values
# (e.g. "8.030" or "1.005") into numeric values, which would silently
# drop meaningful trailing digits / change the cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''29.293.81'
$ws.Range('E2').Value = '''  -0.37%  '
$ws.Range('D3').Value = '''1.832.06'
$ws.Range('E3').Value = '''  -0.51%  '
$ws.Range('D4').Value = '''1.005'
$ws.Range('E4').Value = '''  +0.61%  '
$ws.Range('D5').Value = '''235.55'
$ws.Range('E5').Value = '''  -1.58%  '
$ws.Range('D6').Value = '''0.6029'
$ws.Range('E6').Value = '''  -3.89%  '
$ws.Range('E7').Value = '''  +0.43%  '
$ws.Range('D8').Value = '''0.07063'
$ws.Range('E8').Value = '''  -5.14%  '
$ws.Range('D9').Value = '''0.2802'
$ws.Range('E9').Value = '''  -3.34%  '
$ws.Range('D10').Value = '''23.56'
$ws.Range('E10').Value = '''  -5.12%  '
$ws.Range('D11').Value = '''0.07668'
$ws.Range('E11').Value = '''  -0.59%  '
$ws.Range('D12').Value = '''1.829.33'
$ws.Range('E12').Value = '''  -0.64%  '
$ws.Range('D13').Value = '''4.806'
$ws.Range('E13').Value = '''  -3.40%  '
$ws.Range('D14').Value = '''0.000009944'
$ws.Range('E14').Value = '''  -3.26%  '
$ws.Range('D15').Value = '''0.6289'
$ws.Range('E15').Value = '''  -6.97%  '
$ws.Range('B16').Value = '''Litecoin'
$ws.Range('C16').Value = '''https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').Value = '''79.35'
$ws.Range('E16').Value = '''  -3.04%  '
$ws.Range('B17').Value = '''Uniswap'
$ws.Range('C17').Value = '''https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').Value = '''5.856'
$ws.Range('E17').Value = '''  -6.34%  '
$ws.Range('B18').Value = '''WrappedBTC'
$ws.Range('C18').Value = '''https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '''29.277.34'
$ws.Range('E18').Value = '''  -0.38%  '
$ws.Range('B19').Value = '''BitcoinCash'
$ws.Range('C19').Value = '''https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').Value = '''226.09'
$ws.Range('E19').Value = '''  -3.09%  '
$ws.Range('B20').Value = '''Dai'
$ws.Range('C20').Value = '''https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').Value = '''1.005'
$ws.Range('E20').Value = '''  +0.46%  '
$ws.Range('B21').Value = '''Avalanche'
$ws.Range('C21').Value = '''https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').Value = '''11.72'
$ws.Range('E21').Value = '''  -4.78%  '
$ws.Range('B22').Value = '''Chainlink'
$ws.Range('C22').Value = '''https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D22').Value = '''7.013'
$ws.Range('E22').Value = '''  -4.36%  '
$ws.Range('B23').Value = '''BinanceUSD'
$ws.Range('C23').Value = '''https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D23').Value = '''1.005'
$ws.Range('E23').Value = '''  +0.45%  '
$ws.Range('B24').Value = '''Monero'
$ws.Range('C24').Value = '''https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D24').Value = '''156.46'
$ws.Range('E24').Value = '''  -1.14%  '
$ws.Range('B25').Value = '''Cosmos'
$ws.Range('C25').Value = '''https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D25').Value = '''8.030'
$ws.Range('E25').Value = '''  -5.45%  '
$ws.Range('B26').Value = '''Stellar'
$ws.Range('C26').Value = '''https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D26').Value = '''0.1304'
$ws.Range('E26').Value = '''  -3.48%  '
$ws.Range('B27').Value = '''EthereumClassic'
$ws.Range('C27').Value = '''https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '''16.58'
$ws.Range('E27').Value = '''  -4.52%  '
$ws.Range('B28').Value = '''Toncoin'
$ws.Range('C28').Value = '''https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = '''1.482'
$ws.Range('E28').Value = '''  +1.33%  '
$ws.Range('B29').Value = '''Hedera'
$ws.Range('C29').Value = '''https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D29').Value = '''0.06214'
$ws.Range('E29').Value = '''  -13.19%  '
$ws.Range('B30').Value = '''PancakeSwap'
$ws.Range('C30').Value = '''https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = '''1.446'
$ws.Range('E30').Value = '''  -2.00%  '
$ws.Range('B31').Value = '''Filecoin'
$ws.Range('C31').Value = '''https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '''3.841'
$ws.Range('E31').Value = '''  -4.89%  '
$ws.Range('B32').Value = '''InternetComputer(DFINITY)'
$ws.Range('C32').Value = '''https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '''3.804'
$ws.Range('E32').Value = '''  -6.27%  '
$ws.Range('B33').Value = '''ARBITRUM'
$ws.Range('C33').Value = '''https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').Value = '''1.125'
$ws.Range('E33').Value = '''  -1.28%  '
$ws.Range('B34').Value = '''LidoDAOToken'
$ws.Range('C34').Value = '''https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').Value = '''1.747'
$ws.Range('E34').Value = '''  -3.84%  '
$ws.Range('B35').Value = '''ImmutableX'
$ws.Range('C35').Value = '''https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value = '''0.6437'
$ws.Range('E35').Value = '''  -7.72%  '
$ws.Range('B36').Value = '''HuobiToken'
$ws.Range('C36').Value = '''https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = '''2.543'
$ws.Range('E36').Value = '''  -1.14%  '
$ws.Range('B37').Value = '''Maker'
$ws.Range('C37').Value = '''https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').Value = '''1.221.67'
$ws.Range('E37').Value = '''  -1.15%  '
$ws.Range('B38').Value = '''MXToken'
$ws.Range('C38').Value = '''https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = '''2.737'
$ws.Range('E38').Value = '''  -2.75%  '
$ws.Range('B39').Value = '''VeChain'
$ws.Range('C39').Value = '''https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '''0.01746'
$ws.Range('E39').Value = '''  -5.02%  '
$ws.Range('B40').Value = '''FraxShare'
$ws.Range('C40').Value = '''https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = '''6.564'
$ws.Range('E40').Value = '''  -6.02%  '
$ws.Range('B41').Value = '''TrustWalletToken'
$ws.Range('C41').Value = '''https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = '''0.9054'
$ws.Range('E41').Value = '''  -4.68%  '
$ws.Range('B42').Value = '''PaxDollar'
$ws.Range('C42').Value = '''https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').Value = '''1.004'
$ws.Range('E42').Value = '''  +0.42%  '
$ws.Range('B43').Value = '''RocketPoolETH'
$ws.Range('C43').Value = '''https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D43').Value = '''1.991.67'
$ws.Range('E43').Value = '''  +0.01%  '
$ws.Range('B44').Value = '''Quant'
$ws.Range('C44').Value = '''https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').Value = '''100.85'
$ws.Range('E44').Value = '''  +0.00%  '
$ws.Range('B45').Value = '''Aave'
$ws.Range('C45').Value = '''https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = '''62.69'
$ws.Range('E45').Value = '''  -4.28%  '
$ws.Range('B46').Value = '''BabyDogeCoin'
$ws.Range('C46').Value = '''https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '''0.00000000118'
$ws.Range('E46').Value = '''  -1.72%  '
$ws.Range('B47').Value = '''EnergySwap'
$ws.Range('C47').Value = '''https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '''8.517'
$ws.Range('E47').Value = '''  -4.81%  '
$ws.Range('B48').Value = '''RenderToken'
$ws.Range('C48').Value = '''https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').Value = '''1.585'
$ws.Range('E48').Value = '''  -8.22%  '
$ws.Range('B49').Value = '''Mantle'
$ws.Range('C49').Value = '''https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').Value = '''0.4573'
$ws.Range('E49').Value = '''  -0.24%  '
$ws.Range('B50').Value = '''Cronos'
$ws.Range('C50').Value = '''https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '''0.05517'
$ws.Range('E50').Value = '''  -2.46%  '
$ws.Range('B51').Value = '''Aptos'
$ws.Range('C51').Value = '''https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D51').Value = '''6.452'
$ws.Range('E51').Value = '''  -7.49%  '
